$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 2 (new values)
$ws.Range("A2").Value = 112181853
$ws.Range("B2").Value = 78228
$ws.Range("E2").Value = 6453
$ws.Range("F2").Value = "Vedskivlav"
$ws.Range("G2").Value = "Hertelidea botryosa"
$ws.Range("H2").Value = "(Fr.) Printzen & Kantvilas"
$ws.Range("Q2").Value = 431106
$ws.Range("R2").Value = 6811802

# Row 3 (new values)
$ws.Range("A3").Value = 112182534
$ws.Range("B3").Value = 77388
$ws.Range("E3").Value = 6446
$ws.Range("F3").Value = "Kolflarnlav"
$ws.Range("G3").Value = "Carbonicola anthracophila"
$ws.Range("H3").Value = "(Nyl.) Bendiksby & Timdal"
$ws.Range("Q3").Value = 431104
$ws.Range("R3").Value = 6811805

# Row 4 (new values)
$ws.Range("A4").Value = 112181898
$ws.Range("B4").Value = 78202
$ws.Range("E4").Value = 229821
$ws.Range("F4").Value = "Vedflamlav"
$ws.Range("G4").Value = "Ramboldia elabens"
$ws.Range("H4").Value = "(Fr.) Kantvilas & Elix"
$ws.Range("Q4").Value = 431104
$ws.Range("R4").Value = 6811804
